# Apply "Raw and Clean Data from SSA for June 18th" update:
#  - B18's number format changes from the short "YYYY-MM-DD" style to the
#    long "YYYY-MM-DD HH:MM:SS" style used by the rest of column B.
#  - A new data row (row 19) is appended for June 18, 2020 (serial 44000).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix B18's number format to match the rest of column B (style used by B2:B17).
$ws.Range("B18").NumberFormat = $ws.Range("B17").NumberFormat

# 2) Append the new row of data for June 18th.
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = 44000
$ws.Range("C19").Value = 165455
$ws.Range("D19").Value = 228248
$ws.Range("E19").Value = 59778
$ws.Range("F19").Value = 19747
$ws.Range("G19").Value = 31.86

# Give the new row the same look & feel as the row above it (row 18):
#  - column A keeps the bold/bordered/centered numbering style
#  - column B gets the original short date format ("YYYY-MM-DD") that B18 used
#    before this edit
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("B19").NumberFormat = "YYYY-MM-DD"
